$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark currently sitting at the end of
#    the "...периферийные элементы (вибродвигатель, кнопка и прочее)."
#    paragraph (it is being relocated into the new table row below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the table row whose text contains "переходе на другую
#    частоту" (the "Setup channel" row) and insert two new, blank
#    rows right after it.
# ------------------------------------------------------------------
$targetTable = $null
$targetRowIndex = 0

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
        $rowText = $tbl.Rows.Item($ri).Range.Text
        if ($rowText -like "*переходе на другую частоту*") {
            $targetTable = $tbl
            $targetRowIndex = $ri
        }
    }
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Clear-CellContent($cell, [string]$innerXml) {
    $para = $cell.Range.Paragraphs.Item(1)
    $insPoint = $para.Range
    $insPoint.Collapse(1)
    [void]$insPoint.InsertXML($innerXml)
}

# Plain empty-cell paragraph, matching the rest of the table
# (<w:ind w:firstLine="0"/> only).
$plainParaXml = '<w:p xmlns:w="' + $w + '"><w:pPr><w:ind w:firstLine="0"/></w:pPr></w:p>'

# First new (fully blank) row.
$refRow1 = $targetTable.Rows.Item($targetRowIndex + 1)
$newRow1 = $targetTable.Rows.Add($refRow1)
Clear-CellContent $newRow1.Cells.Item(1) $plainParaXml
Clear-CellContent $newRow1.Cells.Item(2) $plainParaXml
Clear-CellContent $newRow1.Cells.Item(3) $plainParaXml

# Second new row: first cell carries the relocated _GoBack bookmark.
$bookmarkParaXml = '<w:p xmlns:w="' + $w + '"><w:pPr><w:ind w:firstLine="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$refRow2 = $targetTable.Rows.Item($targetRowIndex + 2)
$newRow2 = $targetTable.Rows.Add($refRow2)
Clear-CellContent $newRow2.Cells.Item(1) $bookmarkParaXml
Clear-CellContent $newRow2.Cells.Item(2) $plainParaXml
Clear-CellContent $newRow2.Cells.Item(3) $plainParaXml

Write-Output "Rows inserted after row $targetRowIndex; new row count = $($targetTable.Rows.Count)"
